$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space - matches the separator already used between first
# and last names throughout the existing "name" column data.
$nbsp = [char]0x00A0

# Insert 9 new rows (22-30), one at a time, so each new row inherits the
# cell formatting (styles) of the row immediately above it - matching the
# existing style pattern used by every data row (D: s="2", I: s="1").
for ($i = 22; $i -le 30; $i++) {
    $ws.Rows($i).Insert()
}

# New user_detail rows, following the same layout as the existing data:
# A=id  B=uin  C=name  D=email  E=mobile  F=status_code  G=lang_code
# H=last_login_method  I=is_active  J=cr_by  K=cr_dtimes  L=eff_dtimes
$rows = @(
    @{ r = 22; id = 110021; uin = 7316931025; name = "Magdalena${nbsp}Weber";    email = "magdalena.weber@xyz.com";    mobile = 932122450 },
    @{ r = 23; id = 110022; uin = 9137847236; name = "Adrienne${nbsp}Hoffman";   email = "adrienne.hoffman@xyz.com";   mobile = 848488000 },
    @{ r = 24; id = 110023; uin = 8428758532; name = "Adrienne${nbsp}Mcgee";     email = "adrienne.mcgee@xyz.com";     mobile = 894773246 },
    @{ r = 25; id = 110024; uin = 9804209494; name = "Amare${nbsp}Coleman";      email = "amare.coleman@xyz.com";      mobile = 956554588 },
    @{ r = 26; id = 110025; uin = 7105248214; name = "Dawson${nbsp}Ibarra";      email = "dawson.ibarra@xyz.com";      mobile = 765455583 },
    @{ r = 27; id = 110026; uin = 9316557128; name = "Elvis${nbsp}Mcmillan";     email = "elvis.mcmillan@xyz.com";     mobile = 884282274 },
    @{ r = 28; id = 110027; uin = 8103486949; name = "Steve${nbsp}George";       email = "steve.george@xyz.com";       mobile = 971073663 },
    @{ r = 29; id = 110028; uin = 9601932866; name = "Colton${nbsp}Elliott";     email = "colton.elliott@xyz.com";     mobile = 809908673 },
    @{ r = 30; id = 110029; uin = 9317596765; name = "Carolyn${nbsp}Rodriguez"; email = "carolyn.rodriguez@xyz.com"; mobile = 818876429 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.id
    $ws.Range("B$r").Value = $row.uin
    $ws.Range("C$r").Value = $row.name
    $ws.Range("D$r").Value = $row.email
    $ws.Range("E$r").Value = $row.mobile
    $ws.Range("F$r").Value = "ACT"
    $ws.Range("G$r").Value = "eng"
    $ws.Range("H$r").Value = "PWD"
    $ws.Range("I$r").Value = $true
    $ws.Range("J$r").Value = "superadmin"
    $ws.Range("K$r").Value = "now()"
    $ws.Range("L$r").Value = "now()"
}

# Restore the view: select the newly added id column range for the last
# block of rows (matches the selection left behind after entering this data).
$null = $ws.Range("A22:A30").Select()
